# "break out stock.yaml completed"
# - fix D207:D212 (bsecode) on the "day" sheet: they were written as text,
#   convert them to real numbers (value unchanged).
# - append 9 new rows (213-221) of freshly scraped stock data to the
#   "day" sheet, extending the used range from A1:I212 to A1:I221.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. D207:D212 were stored as text ("541154" etc.) -> store as numbers.
$bseFix = @{
    207 = 541154
    208 = 500247
    209 = 500085
    210 = 513599
    211 = 500103
    212 = 532234
}
foreach ($r in $bseFix.Keys) {
    $ws.Cells.Item($r, 4).Value = $bseFix[$r]
}

# --- 2. Append the new rows scraped on 25/07/2024.
$newRows = @(
    @(213, 1, "MARUTI",     "Maruti Suzuki India Limited",             "532500", 0.18,  12509.2,  396789,    "day", "25/07/2024 11:34:43"),
    @(214, 2, "NESTLEIND",  "Nestle India Limited",                    "500790", -2.39, 2480.65,  3469181,   "day", "25/07/2024 11:34:43"),
    @(215, 3, "DLF",        "Dlf Limited",                             "532868", -0.94, 811.7,    3694096,   "day", "25/07/2024 11:34:43"),
    @(216, 4, "HINDALCO",   "Hindalco Industries Limited",             "500440", -0.78, 646.55,   6172164,   "day", "25/07/2024 11:34:43"),
    @(217, 5, "AUBANK",     "AU Small Finance Bank",                   "540611", -4.23, 631.65,   2402621,   "day", "25/07/2024 11:34:43"),
    @(218, 6, "LAURUSLABS", "Laurus Labs Limited",                     "540222", -0.16, 433.85,   6274677,   "day", "25/07/2024 11:34:43"),
    @(219, 7, "HINDPETRO",  "Hindustan Petroleum Corporation Limited", "500104", 5.62,  373.75,   20611250,  "day", "25/07/2024 11:34:43"),
    @(220, 8, "BPCL",       "Bharat Petroleum Corporation Limited",    "500547", 3.56,  326.15,   30554549,  "day", "25/07/2024 11:34:43"),
    @(221, 9, "BEL",        "Bharat Electronics Limited",              "500049", 0.45,  301.45,   27820556,  "day", "25/07/2024 11:34:43")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]          # A sr
    $ws.Cells.Item($r, 2).Value = $row[2]          # B nsecode
    $ws.Cells.Item($r, 3).Value = $row[3]          # C name
    # D bsecode is written as text by the scraper, even though it looks
    # numeric (mirrors the raw feed -- leading apostrophe forces text so
    # Excel doesn't silently coerce it to a number).
    $ws.Cells.Item($r, 4).Value = "'" + $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]          # E per_chg
    $ws.Cells.Item($r, 6).Value = $row[6]          # F close
    $ws.Cells.Item($r, 7).Value = $row[7]          # G volume
    $ws.Cells.Item($r, 8).Value = $row[8]          # H timeframe
    $ws.Cells.Item($r, 9).Value = $row[9]          # I Date Time
}
